# Remove file upload functionality
# Appends the new capture-log row (row 90) to each of the 4 sheets, mirroring
# the shape of the previous rows. This is purely additive logging data; no
# existing rows/cells are modified.

function Add-LogRow($ws, $row, $timeSerial, $totalLenHex, $idHex, $actualLenHex, $checksumHex, $totalLenDec, $idDec, $actualLenDec, $checksumDec) {
    $ws.Cells.Item($row, 1).Value = $timeSerial
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 2).Value = $totalLenHex
    $ws.Cells.Item($row, 3).Value = $idHex
    $ws.Cells.Item($row, 4).Value = $actualLenHex
    $ws.Cells.Item($row, 5).Value = $checksumHex
    $ws.Cells.Item($row, 6).Value = $totalLenDec
    $ws.Cells.Item($row, 7).Value = $idDec
    $ws.Cells.Item($row, 8).Value = $actualLenDec
    $ws.Cells.Item($row, 9).Value = $checksumDec
}

$wb = $excel.ActiveWorkbook

# --- FE_LFT_#1 ---
$ws1 = $wb.Worksheets.Item("FE_LFT_#1")
$ws1IdDec = [double]"7.598631275147109e+23"
Add-LogRow $ws1 90 45876.49712962963 "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x0C" "0xf" 380 $ws1IdDec 268 15

# --- FE_LFT_#2 ---
$ws2 = $wb.Worksheets.Item("FE_LFT_#2")
$ws2IdDec = [double]"5.68432987514711e+23"
Add-LogRow $ws2 90 45876.49712962963 "0x01,0x90" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x1C" "0xe" 400 $ws2IdDec 284 14

# --- FE_PLT_#1 ---
$ws3 = $wb.Worksheets.Item("FE_PLT_#1")
$ws3IdDec = [double]"5.68631262647114e+23"
Add-LogRow $ws3 90 45876.49712962963 "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x5D" "0x3" 110 $ws3IdDec 93 3

# --- FE_PLT_#2 ---
$ws4 = $wb.Worksheets.Item("FE_PLT_#2")
$ws4IdDec = [double]"9.85046333984776e+23"
Add-LogRow $ws4 90 45876.49712962963 "0x00,0x6e" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x5C" "0x3" 110 $ws4IdDec 92 3
